$d = $word.ActiveDocument
$d.Content.Find.Execute("2024-11-09 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-11-10 Sunday", 2) | Out-Null
$d.Content.Find.Execute("62÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "46÷3=", 2) | Out-Null
$d.Content.Find.Execute("75÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "96÷6=", 2) | Out-Null
$d.Content.Find.Execute("16÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "73÷2=", 2) | Out-Null
$d.Content.Find.Execute("24÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "18÷5=", 2) | Out-Null
$d.Content.Find.Execute("64÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "98÷5=", 2) | Out-Null
$d.Content.Find.Execute("83÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "13÷2=", 2) | Out-Null
$d.Content.Find.Execute("56÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "74÷4=", 2) | Out-Null
$d.Content.Find.Execute("23÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "32÷7=", 2) | Out-Null
$d.Content.Find.Execute("18÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "17÷7=", 2) | Out-Null
$d.Content.Find.Execute("71÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "37÷2=", 2) | Out-Null
$d.Content.Find.Execute("50÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "73÷2=", 2) | Out-Null
$d.Content.Find.Execute("85÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "41÷8=", 2) | Out-Null
$d.Content.Find.Execute("44÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "35÷7=", 2) | Out-Null
$d.Content.Find.Execute("88÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "15÷6=", 2) | Out-Null
$d.Content.Find.Execute("51÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "86÷7=", 2) | Out-Null
$d.Content.Find.Execute("77÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "99÷5=", 2) | Out-Null
$d.Content.Find.Execute("98÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "14÷4=", 2) | Out-Null
$d.Content.Find.Execute("99÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "83÷7=", 2) | Out-Null
$d.Content.Find.Execute("28÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "20÷8=", 2) | Out-Null
$d.Content.Find.Execute("95÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "37÷7=", 2) | Out-Null
$d.Content.Find.Execute("86÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "44÷6=", 2) | Out-Null
$d.Content.Find.Execute("13÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "33÷7=", 2) | Out-Null
$d.Content.Find.Execute("97÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "70÷7=", 2) | Out-Null
$d.Content.Find.Execute("84÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "14÷2=", 2) | Out-Null
$d.Content.Find.Execute("93÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "15÷3=", 2) | Out-Null
